$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, derived from the authoritative diff.
$updates = @{
    'D2' = '30.340.14'
    'E2' = '  -1.20%  '
    'D3' = '1.888.21'
    'E3' = '  -1.64%  '
    'E4' = '  +0.07%  '
    'D5' = '237.90'
    'E5' = '  -1.21%  '
    'D6' = '1.001'
    'E6' = '  +0.05%  '
    'D7' = '0.4816'
    'E7' = '  -2.44%  '
    'D8' = '0.2896'
    'E8' = '  -3.88%  '
    'D9' = '0.06607'
    'E9' = '  -2.57%  '
    'D10' = '1.906.36'
    'E10' = '  -0.59%  '
    'D11' = '16.92'
    'E11' = '  -2.07%  '
    'D12' = '0.07382'
    'E12' = '  +0.67%  '
    'D13' = '5.174'
    'E13' = '  -1.15%  '
    'D14' = '87.79'
    'E14' = '  -1.05%  '
    'D15' = '0.6629'
    'E15' = '  -1.95%  '
    'D16' = '30.326.44'
    'E17' = '  -0.87%  '
    'E18' = '  -0.05%  '
    'D19' = '0.000007756'
    'E19' = '  -2.88%  '
    'D20' = '5.459'
    'E20' = '  +1.20%  '
    'D21' = '2.148.28'
    'E21' = '  -0.61%  '
    'D22' = '1.000'
    'E22' = '  +0.00%  '
    'D23' = '191.67'
    'E23' = '  -3.55%  '
    'D24' = '6.206'
    'E24' = '  -2.13%  '
    'D25' = '9.452'
    'E25' = '  -2.43%  '
    'D26' = '165.36'
    'E26' = '  +1.81%  '
    'D27' = '18.25'
    'E27' = '  -2.26%  '
    'E28' = '  -1.30%  '
    'D29' = '1.452'
    'E29' = '  -1.55%  '
    'E30' = '  -2.58%  '
    'D31' = '0.09174'
    'E31' = '  +0.04%  '
    'D32' = '4.055'
    'E32' = '  -0.78%  '
    'D33' = '0.05078'
    'E33' = '  -3.67%  '
    'D34' = '0.7398'
    'E34' = '  -0.65%  '
    'D35' = '1.144'
    'E35' = '  +1.83%  '
    'E36' = '  +0.25%  '
    'E37' = '  -1.51%  '
    'D38' = '2.648'
    'E38' = '  -2.84%  '
    'D39' = '0.9205'
    'E39' = '  -0.79%  '
    'D40' = '2.090'
    'E40' = '  -0.08%  '
    'D41' = '5.910'
    'E41' = '  -0.71%  '
    'B42' = 'Quant'
    'C42' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D42' = '106.23'
    'E42' = '  -0.82%  '
    'B43' = 'TheSandbox'
    'C43' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D43' = '0.4334'
    'E43' = '  -3.88%  '
    'E44' = '  -0.11%  '
    'D45' = '0.1374'
    'E45' = '  -2.36%  '
    'D46' = '7.654'
    'E46' = '  -0.94%  '
    'E47' = '  +8.83%  '
    'D48' = '65.05'
    'E48' = '  -9.89%  '
    'D49' = '9.011'
    'E49' = '  -0.61%  '
    'D50' = '34.30'
    'E50' = '  -2.86%  '
    'E51' = '  -2.00%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
